# Add two new rows (65 and 66) of data to the bottom of the log on Sheet1,
# matching the date-formatted style already used in column A (A2:A64).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style (date number format) of the last existing data row (A64)
# down into the two new cells, then fill in the actual values. This keeps
# the new cells using the same shared cell style (numFmtId 14, "m/d/yyyy")
# instead of creating a brand new style entry.
$ws.Range("A64").Copy()
$ws.Range("A65:A66").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A65").Value = 46045
$ws.Range("B65").Value = 2

$ws.Range("A66").Value = 46048
$ws.Range("B66").Value = 5

# Match the new selection reflected in the saved workbook.
$ws.Range("A65:B66").Select()
